# Insert two new data rows at the top of the data block (before old row 391),
# shifting the existing rows 391-449 down to 393-451, then populate the
# two newly inserted rows (391-392) with their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 391. This shifts rows 391:449 down to 393:451.
$ws.Rows("391:392").Insert()

# Common (unchanged) values shared by every row in this block.
$mercado = "Agrícola del Norte S.A. de Arica"
$region = "Arica y Parinacota"
$categoriaId = 100112043
$categoria = "Pepino ensalada"
$variedad = "Sin especificar"
$origen = "Región de Arica y Parinacota"
$clasificacion = "Hortaliza"

# New row 391 data
$ws.Cells.Item(391, 1).Value = 1
$ws.Cells.Item(391, 2).Value = $mercado
$ws.Cells.Item(391, 3).Value = $region
$ws.Cells.Item(391, 4).Value = 44995
$ws.Cells.Item(391, 5).Value = 15
$ws.Cells.Item(391, 6).Value = $categoriaId
$ws.Cells.Item(391, 7).Value = $categoria
$ws.Cells.Item(391, 8).Value = $variedad
$ws.Cells.Item(391, 9).Value = "Primera"
$ws.Cells.Item(391, 10).Value = 150
$ws.Cells.Item(391, 11).Value = 7000
$ws.Cells.Item(391, 12).Value = 7500
$ws.Cells.Item(391, 13).Value = 7250
$ws.Cells.Item(391, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(391, 15).Value = $origen
$ws.Cells.Item(391, 16).Value = 104
$ws.Cells.Item(391, 17).Value = 70
$ws.Cells.Item(391, 18).Value = $clasificacion

# New row 392 data
$ws.Cells.Item(392, 1).Value = 1
$ws.Cells.Item(392, 2).Value = $mercado
$ws.Cells.Item(392, 3).Value = $region
$ws.Cells.Item(392, 4).Value = 44995
$ws.Cells.Item(392, 5).Value = 15
$ws.Cells.Item(392, 6).Value = $categoriaId
$ws.Cells.Item(392, 7).Value = $categoria
$ws.Cells.Item(392, 8).Value = $variedad
$ws.Cells.Item(392, 9).Value = "Segunda"
$ws.Cells.Item(392, 10).Value = 160
$ws.Cells.Item(392, 11).Value = 6000
$ws.Cells.Item(392, 12).Value = 6500
$ws.Cells.Item(392, 13).Value = 6250
$ws.Cells.Item(392, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(392, 15).Value = $origen
$ws.Cells.Item(392, 16).Value = 62
$ws.Cells.Item(392, 17).Value = 100
$ws.Cells.Item(392, 18).Value = $clasificacion

# Match the date cell style/format used by the rest of column D (custom
# "YYYY-MM-DD HH:MM:SS" number format), same as the cells that shifted down.
$ws.Cells.Item(391, 4).NumberFormat = $ws.Cells.Item(393, 4).NumberFormat
$ws.Cells.Item(392, 4).NumberFormat = $ws.Cells.Item(393, 4).NumberFormat
